# edit.ps1 - PowerPoint COM-interop script (PowerShell-style) applying the
# changes captured by the target diff:
#
#   1. Three tables (on slides 14, 15 and 16) get their table style switched
#      from {749B339B-E2E6-45D4-8417-27845DA91901} to
#      {3CCFC570-4BFB-42BD-81AD-724683C0BD1E}.
#   2. The deck's theme ("Integral" / "Red Violet" colour scheme) is swapped
#      for the default "Office Theme" colour scheme (the twelve theme
#      colours are updated to the standard Office palette).

$p = $ppt.ActivePresentation

# --- 1. Update the table styles on the three affected slides -------------
$oldStyleId = "{749B339B-E2E6-45D4-8417-27845DA91901}"
$newStyleId = "{3CCFC570-4BFB-42BD-81AD-724683C0BD1E}"

$tableSlideIndexes = @(14, 15, 16)
foreach ($slideIndex in $tableSlideIndexes) {
    $slide = $p.Slides.Item($slideIndex)
    for ($shapeIndex = 1; $shapeIndex -le $slide.Shapes.Count; $shapeIndex++) {
        $shape = $slide.Shapes.Item($shapeIndex)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}

# --- 2. Re-colour the theme to the standard "Office" palette -------------
function RGBVal($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$slide1 = $p.Slides.Item(1)
$themeColors = $slide1.ThemeColorScheme

$themeColors.Item(1).RGB  = RGBVal 0x00 0x00 0x00   # dk1
$themeColors.Item(2).RGB  = RGBVal 0xFF 0xFF 0xFF   # lt1
$themeColors.Item(3).RGB  = RGBVal 0x44 0x54 0x6A   # dk2
$themeColors.Item(4).RGB  = RGBVal 0xE7 0xE6 0xE6   # lt2
$themeColors.Item(5).RGB  = RGBVal 0x5B 0x9B 0xD5   # accent1
$themeColors.Item(6).RGB  = RGBVal 0xED 0x7D 0x31   # accent2
$themeColors.Item(7).RGB  = RGBVal 0xA5 0xA5 0xA5   # accent3
$themeColors.Item(8).RGB  = RGBVal 0xFF 0xC0 0x00   # accent4
$themeColors.Item(9).RGB  = RGBVal 0x44 0x72 0xC4   # accent5
$themeColors.Item(10).RGB = RGBVal 0x70 0xAD 0x47   # accent6
$themeColors.Item(11).RGB = RGBVal 0x05 0x63 0xC1   # hlink
$themeColors.Item(12).RGB = RGBVal 0x95 0x4F 0x72   # folHlink
